$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new drive entry (row 10) for the 230V, 3A / CIMR-VU2A0004BMA-092 unit.
# Description / manufacturer part number are populated first.
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 21573
$ws.Range("C10").Value = "VARIABLE FREQUENCY DRIVE, 230V, 3A"
$ws.Range("D10").Value = "Yaskawa"
$ws.Range("E10").Value = "CIMR-VU2A0004BMA-092"
$ws.Range("H10").Value = "V1000_GRP_DESC"

# Update PARAM_ND_LIST / PARAM_HD_LIST columns (F, G) for existing rows 2-9
# to use the part-number-based naming convention instead of the generic
# V1000_xV_ND / V1000_xV_HD values.
$ws.Range("F2").Value = "21520_ND"
$ws.Range("G2").Value = "21520_HD"

$ws.Range("F3").Value = "21521_ND"
$ws.Range("G3").Value = "21521_HD"

$ws.Range("F4").Value = "21522_ND"
$ws.Range("G4").Value = "21522_HD"

$ws.Range("F5").Value = "21523_ND"
$ws.Range("G5").Value = "21523_HD"

$ws.Range("F6").Value = "21524_ND"
$ws.Range("G6").Value = "21524_HD"

$ws.Range("F7").Value = "21525_ND"
$ws.Range("G7").Value = "21525_HD"

$ws.Range("F8").Value = "21540_ND"
$ws.Range("G8").Value = "21540_HD"

$ws.Range("F9").Value = "21541_ND"
$ws.Range("G9").Value = "21541_HD"

# Finally, the new row's PARAM_ND_LIST / PARAM_HD_LIST values.
$ws.Range("F10").Value = "21573_ND"
$ws.Range("G10").Value = "21573_HD"

# Move the active selection to reflect where the user left off editing.
$ws.Range("H11").Select()
